$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each record.
# Every data row (2-143) currently shows the same value (45172, i.e.
# 2023-09-03); the update bumps it to 45175 (2023-09-06) for all of them.
$ws.Range("C2:C143").Value = 45175
